$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / sample-size-ish values)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2
$ws.Range("B2").Value = 8.4136337932082483
$ws.Range("C2").Value = 18.795735196116311
$ws.Range("D2").Value = 23.730625264933053
$ws.Range("E2").Value = 23.821750761293913

# Row 3
$ws.Range("B3").Value = 11.526169832176834
$ws.Range("C3").Value = 14.318636710786166
$ws.Range("D3").Value = 12.154062769443755
$ws.Range("E3").Value = 27.355688807826709

# Narrow the active selection to match the edited range
$ws.Range("B1:E3").Select()
